# "Generate Report for Archive"
# The localization status for these two files moved on from handoff: the
# "Ready for handoff" status is now "In Translation" everywhere it is used
# (Overview!E2:F3, and the "Status" column on the per-language sheets).
# Narrowing that now-shorter text also lets the Status columns shrink a bit.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$narrowerWidth = 12.5   # closest ColumnWidth this host can persist to ~13.33 chars

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $narrowerWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowerWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowerWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowerWidth
